$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "This "
$ws.Range("E7").Value = 8

$ws.Range("E7").Select()
